# DMS: Fix KpiItem Export
# The "Loại KPI sản phẩm" row (row 6) showed a broken placeholder that
# referenced a nested object ({{KpiItems.KpiItemType.Name}}) instead of the
# flattened property the report engine actually exposes
# ({{KpiItems.KpiItemTypeName}}). Fix the merged cell E6:F6 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "{{KpiItems.KpiItemTypeName}}"

# Re-select the cell that was actually fixed, matching the cursor position
# left behind in the saved file.
$ws.Range("E6:F6").Select()
